# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, mirroring the
# existing header style used by column H, and fills in the per-row values
# for rows 2 through 79.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: copy style from H1 onto I1/J1, then set values.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row I/J values for rows 2..79 (index 0 => row 2)
$iValues = @(6,5,6,7,8,7,7,8,7,6,10,10,7,8,7,8,7,8,7,9,7,8,9,7,6,6,7,9,8,9,9,6,8,6,8,10,8,8,8,8,9,9,8,8,8,8,10,8,6,8,5,9,8,4,6,9,8,8,8,9,9,9,9,9,9,6,7,8,9,8,6,7,6,5,6,4,4,4)
$jValues = @(6,5,6,7,8,7,7,8,7,6,10,10,7,8,7,8,7,8,7,9,8,8,9,7,6,6,7,9,8,9,9,6,8,6,8,10,8,8,8,8,9,9,8,8,8,8,10,8,7,8,5,9,8,5,6,9,8,8,8,9,9,9,9,9,9,6,7,8,9,8,6,7,6,6,6,4,4,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
